$d = $word.ActiveDocument

# Locate the empty "List Paragraph" entry that currently has no numbering
# and no runs -- it is the blank changelog slot to be filled in with the
# new 0812527 / 27-05-2012 entry and its description.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Style.NameLocal -eq "List Paragraph" -and $cand.Range.Text.Trim().Length -eq 0) {
        $target = $cand
        break
    }
}

$idx = $target.Index
$target.Range.InsertParagraphAfter()

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$headerXml = '<w:p ' + $wns + '>' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="ListParagraph"/>' + `
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
      '<w:b/>' + `
      '<w:sz w:val="24"/><w:szCs w:val="24"/>' + `
      '<w:lang w:val="en-US"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
      '<w:b/>' + `
      '<w:sz w:val="24"/><w:szCs w:val="24"/>' + `
      '<w:lang w:val="en-US"/>' + `
    '</w:rPr>' + `
    '<w:t>0812527 &#8211; 27/05/2012 19:53</w:t>' + `
  '</w:r>' + `
'</w:p>'

$bodyXml = '<w:p ' + $wns + '>' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="ListParagraph"/>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
      '<w:sz w:val="24"/><w:szCs w:val="24"/>' + `
      '<w:lang w:val="en-US"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
      '<w:sz w:val="24"/><w:szCs w:val="24"/>' + `
    '</w:rPr>' + `
    '<w:t>Trong store Đăng ký đồ án, tui thấy có declare một biến là ThoiHanNop rồi sau đó có so sánh nhưng trước đó không có set gì, vì vậy tui đã thêm lệnh set vào trước câu select.</w:t>' + `
  '</w:r>' + `
'</w:p>'

$null = $d.Paragraphs.Item($idx).Range.InsertXML($headerXml)
$null = $d.Paragraphs.Item($idx + 1).Range.InsertXML($bodyXml)
